$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Reword "You don't need to pass any data between the views." into
#    "You only need to pass data from one view to one of the other views."
#    The final text ends up split across five runs (matching how Word leaves
#    things after a real, incremental editing session):
#      - " to switch between at least three views. You "   (pre-existing run)
#      - "only need to pass data from one view to one of t" (new run)
#      - "he"                                                (new run)
#      - " other"                                            (new run)
#      - " views. "                                          (new run)
# ---------------------------------------------------------------------------
$oldSentence = "You don" + [char]0x2019 + "t need to pass any data between the views. "
$newSentence = "You only need to pass data from one view to one of the other views. "

$target = $d.Content
$target.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $target.Start
$end = $target.End

# Replace the whole sentence in one shot first -- this keeps the edit a single
# (rsid-less) run for now; we carve it into the final run boundaries next.
$whole = $d.Range($start, $end)
$whole.Text = $newSentence
$wholeEnd = $whole.End

$piece2 = "only need to pass data from one view to one of t"
$piece3 = "he"
$piece4 = " other"

$splitAfterYou    = $start + ("You ").Length
$splitAfterPiece2 = $splitAfterYou + $piece2.Length
$splitAfterPiece3 = $splitAfterPiece2 + $piece3.Length
$splitAfterPiece4 = $splitAfterPiece3 + $piece4.Length

foreach ($p in @($splitAfterYou, $splitAfterPiece2, $splitAfterPiece3, $splitAfterPiece4)) {
    $s = $d.Range($p, $wholeEnd)
    $s.Bold = 1
    $s.Bold = 0
}

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark: delete it from the trailing empty paragraph
#    and re-create it right after the sentence we just edited.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($wholeEnd, $wholeEnd))
